$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has two copies of the same "line amp" table, separated by a
# 2-row gap (rows 8-9 are blank). Delete the extra blank row 8 so the
# second table (previously rows 10-14) shifts up to rows 9-13, leaving
# just a single blank row between the tables.
$ws.Rows("8:8").Delete()

# In the first table, the "8" row's 75 kVA/Random column should read
# "8-0_12.9kW" (matching the other 75 kVA column), not "8-8_12.9kW".
$ws.Range("D7").Value = "8-0_12.9kW"

# In the second table (now shifted up), add the missing "4-0_12.9kW"
# data point for the 50 kVA/Set column of the "4" row.
$ws.Range("D12").Value = "4-0_12.9kW"

# Restore the saved view state: zoomed to 85%, scrolled down a bit, with
# D12 as the active selection.
$win = $excel.ActiveWindow
$win.Zoom = 85
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("D12").Select()
